$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORT_DATE (stored as plain text, not a date value)
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# Numeric columns updated per diff
$ws.Range("O2").Value = 2100877375.3
$ws.Range("P2").Value = 294.8877841047
$ws.Range("Q2").Value = 45921946970.85
$ws.Range("R2").Value = 6445.7932401084
$ws.Range("S2").Value = 1945121654.06
$ws.Range("T2").Value = 273.0252708337
$ws.Range("U2").Value = -1948450605.78
$ws.Range("V2").Value = -273.4925361809
$ws.Range("Y2").Value = 1951892252.12
$ws.Range("Z2").Value = 273.9756198081
$ws.Range("AA2").Value = -873547811.54
$ws.Range("AB2").Value = -122.6147615673
$ws.Range("AC2").Value = -712432826.5
$ws.Range("AD2").Value = -155.5861797275
